# Updates the "exp" (column B) placeholder values on the four level-curve
# sheets (tough_levels, intelligent_levels, wise_levels, intuitive_levels).
#
#   tough_levels:        row3 B (exp) 16 -> 8
#   intelligent_levels:  row2 B (exp) 2  -> 4 ; row3 B (exp) 16 -> 12
#   wise_levels:         row2 B (exp) 2  -> 4 ; row3 B (exp) 16 -> 12
#   intuitive_levels:    row2 B (exp) 2  -> 4 ; row3 B (exp) 16 -> 12

$wb = $excel.ActiveWorkbook

$toughLevels = $wb.Worksheets.Item("tough_levels")
$toughLevels.Range("B3").Value = 8

$intelligentLevels = $wb.Worksheets.Item("intelligent_levels")
$intelligentLevels.Range("B2").Value = 4
$intelligentLevels.Range("B3").Value = 12

$wiseLevels = $wb.Worksheets.Item("wise_levels")
$wiseLevels.Range("B2").Value = 4
$wiseLevels.Range("B3").Value = 12

$intuitiveLevels = $wb.Worksheets.Item("intuitive_levels")
$intuitiveLevels.Range("B2").Value = 4
$intuitiveLevels.Range("B3").Value = 12

# Restore the per-sheet cursor/selection state recorded in the edit.
$toughLevels.Activate() | Out-Null
$toughLevels.Range("C11").Select() | Out-Null

$intelligentLevels.Activate() | Out-Null
$intelligentLevels.Range("D17").Select() | Out-Null

$wiseLevels.Activate() | Out-Null
$wiseLevels.Range("G30").Select() | Out-Null

$intuitiveLevels.Activate() | Out-Null
$intuitiveLevels.Range("C11").Select() | Out-Null
